$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '37.791.64'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '2.077.30'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''233.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = '''58.50'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("E11").Value = '  +2.92%  '
$ws.Range("D12").Value = '''14.85'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '2.382.54'
$ws.Range("E13").Value = '  -0.56%  '
$ws.Range("D14").Value = '''21.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.92%  '
$ws.Range("D15").Value = '''0.781'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '''5.35'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").Value = '2.053.62'
$ws.Range("E17").Value = '  -1.47%  '
$ws.Range("D18").Value = '37.688.08'
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("D19").Value = '''6.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.72%  '
$ws.Range("D20").Value = '''71.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = '0.0₃0839'
$ws.Range("E21").Value = '  +1.20%  '
$ws.Range("D22").Value = '''229.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("E24").Value = '  -0.81%  '
$ws.Range("D25").Value = '''2.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").Value = '''9.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.29%  '
$ws.Range("D27").Value = '''172.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("E28").Value = '  -1.29%  '
$ws.Range("D29").Value = '''19.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("E30").Value = '  -1.87%  '
$ws.Range("E31").Value = '  +1.15%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("E34").Value = '  -1.08%  '
$ws.Range("D35").Value = '''2.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.92%  '
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("D37").Value = '''3.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.28%  '
$ws.Range("D38").Value = '''1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("E40").Value = '  +8.51%  '
$ws.Range("D41").Value = '''101.12'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.19%  '
$ws.Range("D42").Value = '''0.0971'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.24%  '
$ws.Range("D43").Value = '''2.93'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("D44").Value = '''16.88'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.83%  '
$ws.Range("D45").Value = '1.446.42'
$ws.Range("E45").Value = '  -0.96%  '
$ws.Range("E46").Value = '  -1.88%  '
$ws.Range("E47").Value = '  -1.14%  '
$ws.Range("D48").Value = '''4.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.15%  '
$ws.Range("E49").Value = '  -1.83%  '
$ws.Range("E50").Value = '  -1.94%  '
$ws.Range("D51").Value = '2.268.98'